$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price/Volume/Hora columns as Text so that numeric-looking
# strings (e.g. "304.01", "6.19%", "15") are stored as literal text, matching
# the source data which stores these as inline strings, not numbers.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = '304.01'
$ws.Range("E2").Value = '6.19%'
$ws.Range("G2").Value = '15'
$ws.Range("E3").Value = '9.27%'
$ws.Range("G3").Value = '15'
$ws.Range("D4").Value = '5.289'
$ws.Range("E4").Value = '4.28%'
$ws.Range("G4").Value = '15'
$ws.Range("D5").Value = '0.07459'
$ws.Range("E5").Value = '10.40%'
$ws.Range("G5").Value = '15'
$ws.Range("D6").Value = '7.848'
$ws.Range("E6").Value = '6.80%'
$ws.Range("G6").Value = '15'
$ws.Range("E7").Value = '9.74%'
$ws.Range("G7").Value = '15'
$ws.Range("D8").Value = '1.470'
$ws.Range("E8").Value = '6.62%'
$ws.Range("G8").Value = '15'
$ws.Range("D9").Value = '0.9164'
$ws.Range("E9").Value = '1.74%'
$ws.Range("G9").Value = '15'
$ws.Range("D10").Value = '0.01734'
$ws.Range("E10").Value = '2,572.51%'
$ws.Range("G10").Value = '15'
$ws.Range("D11").Value = '0.1696'
$ws.Range("E11").Value = '6.42%'
$ws.Range("G11").Value = '15'
$ws.Range("E12").Value = '11.38%'
$ws.Range("G12").Value = '15'
$ws.Range("D13").Value = '0.08039'
$ws.Range("E13").Value = '6.23%'
$ws.Range("G13").Value = '15'
$ws.Range("D14").Value = '0.03029'
$ws.Range("E14").Value = '3.73%'
$ws.Range("G14").Value = '15'
$ws.Range("D15").Value = '0.09904'
$ws.Range("E15").Value = '10.06%'
$ws.Range("G15").Value = '15'
$ws.Range("D16").Value = '0.001492'
$ws.Range("E16").Value = '-4.68%'
$ws.Range("G16").Value = '15'
$ws.Range("D17").Value = '0.04571'
$ws.Range("E17").Value = '1.95%'
$ws.Range("G17").Value = '15'
$ws.Range("D18").Value = '0.006165'
$ws.Range("E18").Value = '-6.09%'
$ws.Range("G18").Value = '15'
$ws.Range("D19").Value = '3.482'
$ws.Range("E19").Value = '0.73%'
$ws.Range("G19").Value = '15'
$ws.Range("D20").Value = '2.229'
$ws.Range("G20").Value = '15'
$ws.Range("D21").Value = '0.3301'
$ws.Range("E21").Value = '2.95%'
$ws.Range("G21").Value = '15'
$ws.Range("D22").Value = '0.1344'
$ws.Range("E22").Value = '1.83%'
$ws.Range("G22").Value = '15'
$ws.Range("D23").Value = '4.523'
$ws.Range("E23").Value = '13.04%'
$ws.Range("G23").Value = '15'
$ws.Range("D24").Value = '0.1624'
$ws.Range("E24").Value = '4.43%'
$ws.Range("G24").Value = '15'
$ws.Range("E25").Value = '1.32%'
$ws.Range("G25").Value = '15'
$ws.Range("D26").Value = '0.004413'
$ws.Range("E26").Value = '0.77%'
$ws.Range("G26").Value = '15'
$ws.Range("D27").Value = '0.0001403'
$ws.Range("E27").Value = '20.18%'
$ws.Range("G27").Value = '15'
$ws.Range("D28").Value = '0.0001779'
$ws.Range("E28").Value = '10.22%'
$ws.Range("G28").Value = '15'
$ws.Range("G29").Value = '15'
$ws.Range("G30").Value = '15'
$ws.Range("G31").Value = '15'
$ws.Range("G32").Value = '15'
$ws.Range("G33").Value = '15'
$ws.Range("G34").Value = '15'
$ws.Range("G35").Value = '15'
$ws.Range("G36").Value = '15'
$ws.Range("G37").Value = '15'
$ws.Range("G38").Value = '15'
$ws.Range("G39").Value = '15'
$ws.Range("D40").Value = '0.04506'
$ws.Range("E40").Value = '6.09%'
$ws.Range("G40").Value = '15'
$ws.Range("D41").Value = '0.007222'
$ws.Range("E41").Value = '6.07%'
$ws.Range("G41").Value = '15'
$ws.Range("D42").Value = '0.1343'
$ws.Range("E42").Value = '8.24%'
$ws.Range("G42").Value = '15'
$ws.Range("D43").Value = '0.002224'
$ws.Range("E43").Value = '1.36%'
$ws.Range("G43").Value = '15'
$ws.Range("D44").Value = '0.01368'
$ws.Range("E44").Value = '20.03%'
$ws.Range("G44").Value = '15'
$ws.Range("D45").Value = '0.00006155'
$ws.Range("E45").Value = '7.53%'
$ws.Range("G45").Value = '15'
$ws.Range("D46").Value = '0.7082'
$ws.Range("E46").Value = '-63.29%'
$ws.Range("G46").Value = '15'
$ws.Range("E47").Value = '-0.11%'
$ws.Range("G47").Value = '15'
$ws.Range("G48").Value = '15'
$ws.Range("G49").Value = '15'
$ws.Range("G50").Value = '15'
$ws.Range("G51").Value = '15'
